$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.395.10'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '2.294.44'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.60'
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.96'
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +2.91%  '
$ws.Range("D9").Value = '2.288.84'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.101'
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.78'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = '2.702.77'
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").Value = '58.329.79'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '2.277.47'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.59'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '315.24'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.01'
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.98'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  -1.97%  '
$ws.Range("D31").Value = '0.0₃0728'
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").Value = '  +2.14%  '
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.87'
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.95'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '291.49'
$ws.Range("E41").Value = '  -2.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.81'
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.46'
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0496'
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.34'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0212'
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("E51").Value = '  +1.14%  '
